$wb = $excel.ActiveWorkbook

# --- Work on the "Repayment schedule" sheet ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (pushes old N/O/P -> O/P/Q)
$ws.Range("N1").EntireColumn.Insert()

# The newly inserted column inherits the neighbouring column's width;
# match the source column's character width as closely as possible.
$ws.Columns.Item(14).ColumnWidth = 9.8

# Update the current selection on this sheet
[void]$ws.Range("S5").Select()

# --- Make "Repayment schedule" the active tab (was "Transactions") ---
# Activating the sheet above already updates the workbook's active tab
# and clears the previous tabSelected flag on "Transactions".
